$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Entries": renumber id column (A4:A8) +1 and move selection
# ---------------------------------------------------------------------
$wsEntries = $wb.Worksheets.Item("Entries")
$wsEntries.Cells.Item(4,1).Value = 1
$wsEntries.Cells.Item(5,1).Value = 2
$wsEntries.Cells.Item(6,1).Value = 3
$wsEntries.Cells.Item(7,1).Value = 4
$wsEntries.Cells.Item(8,1).Value = 5

# ---------------------------------------------------------------------
# Sheet "VatTypes": no data change, only move the selected cell
# ---------------------------------------------------------------------
$wsVat = $wb.Worksheets.Item("VatTypes")

# ---------------------------------------------------------------------
# Sheet "CashBook": insert a new "code" column (string) before the
# existing "transaction_code" column, shifting B:G -> C:H, and fill in
# the new column's metadata/header/data rows.
# ---------------------------------------------------------------------
$wsCash = $wb.Worksheets.Item("CashBook")
$wsCash.Columns.Item(2).Insert()

$wsCash.Cells.Item(1,2).Value = 1
$wsCash.Cells.Item(2,2).Value = "string"
$wsCash.Cells.Item(3,2).Value = "code"
$wsCash.Cells.Item(4,2).Value = 1
$wsCash.Cells.Item(5,2).Value = 2

# ---------------------------------------------------------------------
# Restore the selections on every sheet (must happen after edits so the
# insert on CashBook doesn't disturb the other sheets' selections), and
# keep CashBook as the active/selected tab like in the source file.
# ---------------------------------------------------------------------
$wsEntries.Activate()
$wsEntries.Range("H21").Select()

$wsVat.Activate()
$wsVat.Range("A4").Select()

$wsCash.Activate()
$wsCash.Range("A6").Select()
